# Updated cryptos list on Sat Nov 30 23:58:32 UTC 2024 with GitHub Actions
#
# Refresh live-ish crypto market data (price + 1h volume change) pulled in
# by the daily GitHub Actions job. A couple of rows (ARBITRUM / EnergySwap)
# also swapped rank position this run, so those two rows get their full
# Coin/Link/Price/Volume set rewritten rather than a single-cell bump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Price values (column D) that are plain
# decimals are quote-prefixed so Excel stores them as text (matching the
# sheet's existing inlineStr convention) instead of coercing them to a
# number and silently dropping significant trailing zeros (e.g. "13.40").
$updates = [ordered]@{
    'D2' = '96.487.37'
    'E2' = '  -0.99%  '
    'D3' = '3.708.19'
    'E3' = '  +3.17%  '
    'E4' = '  +0.02%  '
    'D5' = '''237.87'
    'E5' = '  -2.28%  '
    'E6' = '  +8.08%  '
    'D7' = '''654.48'
    'E7' = '  -0.06%  '
    'E8' = '  -0.90%  '
    'E9' = '  +0.35%  '
    'E10' = '  +0.02%  '
    'D11' = '3.707.19'
    'E11' = '  +3.24%  '
    'D12' = '''44.93'
    'E12' = '  +0.14%  '
    'E13' = '  +0.39%  '
    'E14' = '  +5.92%  '
    'D15' = '4.401.27'
    'E15' = '  +3.26%  '
    'E16' = '  +2.87%  '
    'D17' = '96.259.32'
    'E17' = '  -0.90%  '
    'E18' = '  +15.90%  '
    'D19' = '3.720.95'
    'E19' = '  +3.57%  '
    'D20' = '''19.04'
    'E20' = '  +4.31%  '
    'E21' = '  +1.69%  '
    'D22' = '''0.527'
    'E22' = '  -3.63%  '
    'D23' = '''524.53'
    'E23' = '  +1.31%  '
    'E24' = '  +0.08%  '
    'D25' = '''7.01'
    'E25' = '  +0.66%  '
    'D26' = '''0.0000202'
    'E26' = '  -1.14%  '
    'D27' = '''102.81'
    'E27' = '  -0.87%  '
    'D28' = '''13.40'
    'E28' = '  +0.93%  '
    'E29' = '  -5.11%  '
    'D30' = '''12.47'
    'E30' = '  +3.41%  '
    'E31' = '  +1.80%  '
    'E32' = '  +0.07%  '
    'E33' = '  +10.23%  '
    'E34' = '  -1.26%  '
    'D35' = '''676.97'
    'E35' = '  +9.68%  '
    'D36' = '''32.74'
    'E36' = '  +2.66%  '
    'E37' = '  +0.26%  '
    'D38' = '''0.599'
    'E38' = '  +2.72%  '
    'E39' = '  +1.47%  '
    'D40' = '''7.11'
    'E40' = '  +15.81%  '
    'E41' = '  +4.56%  '
    'B42' = 'EnergySwap'
    'C42' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D42' = '''40.27'
    'E42' = '  +23.06%  '
    'B43' = 'ARBITRUM'
    'C43' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D43' = '''0.975'
    'E43' = '  +4.83%  '
    'E44' = '  +2.45%  '
    'D46' = '''0.0459'
    'E46' = '  +1.23%  '
    'D47' = '''0.441'
    'E47' = '  -0.11%  '
    'E48' = '  -1.02%  '
    'D49' = '''23.62'
    'E49' = '  -0.21%  '
    'E50' = '  -1.70%  '
    'D51' = '''3.55'
    'E51' = '  +2.49%  '

}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
